$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.620934247970581
$ws.Range("B1").Value = 1.85272204875946
$ws.Range("C1").Value = 2.310180187225342
$ws.Range("D1").Value = 3.671397924423218
$ws.Range("E1").Value = 2.734396934509277
